# Dashboard page elements visibility, elements heading visibility and
# heading assertion added.
#
# Adds a new "dashboardElements" worksheet (after the existing sheets)
# containing a header row of 7 "Element N Header" columns and a data row
# with the dashboard widget names, then makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the last existing sheet --------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "dashboardElements"

# --- Header row (row 1) -------------------------------------------------
$headers = @(
    "Element 1 Header",
    "Element 2 Header",
    "Element 3 Header",
    "Element 4 Header",
    "Element 5 Header",
    "Element 6 Header",
    "Element 7 Header"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Data row (row 2) ----------------------------------------------------
$data = @(
    "Time at Work",
    "My Actions",
    "Quick Launch",
    "Buzz Latest Posts",
    "Employees on Leave Today",
    "Employee Distribution by Sub Unit",
    "Employee Distribution by Location"
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $newSheet.Cells.Item(2, $i + 1).Value = $data[$i]
}

# --- Formatting: header row bold / centered / wrapped --------------------
$headerRange = $newSheet.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true
$newSheet.Rows.Item(1).RowHeight = 30

# --- Formatting: data row, regular (non-bold) font ------------------------
$dataRange = $newSheet.Range("A2:G2")
$dataRange.Font.Bold = $false
$dataRange.Font.Name = "Calibri"
$dataRange.Font.Size = 11

# --- Column widths (best-fit sized to content) ----------------------------
$newSheet.Columns.Item(1).ColumnWidth = 12.33
$newSheet.Columns.Item(2).ColumnWidth = 10.17
$newSheet.Columns.Item(3).ColumnWidth = 12.33
$newSheet.Columns.Item(4).ColumnWidth = 17.5
$newSheet.Columns.Item(5).ColumnWidth = 24.67
$newSheet.Columns.Item(6).ColumnWidth = 34
$newSheet.Columns.Item(7).ColumnWidth = 34

# --- Page setup ------------------------------------------------------------
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# --- Make the new sheet the active / selected tab ---------------------------
$newSheet.Activate()
$newSheet.Select() | Out-Null
$newSheet.Range("A1").Select() | Out-Null
